# chore: update Sheets via scheduled runner
# Refreshes the market-board derived columns (currentAveragePrice*,
# LevePrice*, LeveProfit*) on each crafting-job sheet with newly
# fetched values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1639.1
$ws.Range("I40").Value = 733.6667
$ws.Range("J40").Value = 2027.1428
$ws.Range("K40").Value = 733.6667
$ws.Range("L40").Value = 2027.1428
$ws.Range("M40").Value = -558.6667
$ws.Range("N40").Value = -2377.1428

$ws.Range("H62").Value = 10482.637
$ws.Range("I62").Value = 1538.625
$ws.Range("J62").Value = 34333.332
$ws.Range("K62").Value = 1538.625
$ws.Range("L62").Value = 34333.332
$ws.Range("M62").Value = -914.625
$ws.Range("N62").Value = -35581.332

$ws.Range("H65").Value = 10482.637
$ws.Range("I65").Value = 1538.625
$ws.Range("J65").Value = 34333.332
$ws.Range("K65").Value = 7693.125
$ws.Range("L65").Value = 171666.66
$ws.Range("M65").Value = -4573.125
$ws.Range("N65").Value = -177906.66

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 38349
$ws.Range("J44").Value = 38349
$ws.Range("L44").Value = 38349
$ws.Range("N44").Value = -39325

$ws.Range("H55").Value = 25806.334
$ws.Range("J55").Value = 25806.334
$ws.Range("L55").Value = 25806.334
$ws.Range("N55").Value = -26436.334

$ws.Range("H74").Value = 574.80554
$ws.Range("I74").Value = 555.087
$ws.Range("J74").Value = 609.6923
$ws.Range("K74").Value = 555.087
$ws.Range("L74").Value = 609.6923
$ws.Range("M74").Value = 318.913
$ws.Range("N74").Value = -2357.6923

$ws.Range("H77").Value = 574.80554
$ws.Range("I77").Value = 555.087
$ws.Range("J77").Value = 609.6923
$ws.Range("K77").Value = 2775.435
$ws.Range("L77").Value = 3048.4615
$ws.Range("M77").Value = 1592.565
$ws.Range("N77").Value = -11784.4615

$ws.Range("H80").Value = 25573.4
$ws.Range("J80").Value = 25573.4
$ws.Range("L80").Value = 25573.4
$ws.Range("N80").Value = -27569.4

$ws.Range("H83").Value = 25573.4
$ws.Range("J83").Value = 25573.4
$ws.Range("L83").Value = 76720.20000000001
$ws.Range("N83").Value = -86704.20000000001

$ws.Range("H122").Value = 2259.5293
$ws.Range("I122").Value = 1926
$ws.Range("J122").Value = 3060
$ws.Range("K122").Value = 5778
$ws.Range("L122").Value = 9180
$ws.Range("M122").Value = -3328
$ws.Range("N122").Value = -14080

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 34371.332
$ws.Range("J35").Value = 34371.332
$ws.Range("L35").Value = 34371.332
$ws.Range("N35").Value = -34991.332

$ws.Range("H60").Value = 39000
$ws.Range("J60").Value = 39000
$ws.Range("L60").Value = 39000
$ws.Range("N60").Value = -40198

$ws.Range("H82").Value = 27608.143
$ws.Range("J82").Value = 37348.6
$ws.Range("L82").Value = 37348.6
$ws.Range("N82").Value = -38114.6

$ws.Range("H85").Value = 27608.143
$ws.Range("J85").Value = 37348.6
$ws.Range("L85").Value = 37348.6
$ws.Range("N85").Value = -40000.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 16365
$ws.Range("J41").Value = 16365
$ws.Range("L41").Value = 16365
$ws.Range("N41").Value = -17221

$ws.Range("H51").Value = 8950
$ws.Range("J51").Value = 8950
$ws.Range("L51").Value = 8950
$ws.Range("N51").Value = -10422

$ws.Range("H55").Value = 6850
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 6850
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 6850
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -7480

$ws.Range("H60").Value = 11062.4
$ws.Range("J60").Value = 11062.4
$ws.Range("L60").Value = 11062.4
$ws.Range("N60").Value = -12084.4

$ws.Range("H61").Value = 8950
$ws.Range("J61").Value = 8950
$ws.Range("L61").Value = 8950
$ws.Range("N61").Value = -9646

$ws.Range("H99").Value = 1728.1034
$ws.Range("I99").Value = 1757.5555
$ws.Range("J99").Value = 1679.909
$ws.Range("K99").Value = 1757.5555
$ws.Range("L99").Value = 1679.909
$ws.Range("M99").Value = -259.5554999999999
$ws.Range("N99").Value = -4675.909

$ws.Range("H109").Value = 32485
$ws.Range("J109").Value = 32485
$ws.Range("L109").Value = 32485
$ws.Range("N109").Value = -34565

$ws.Range("H126").Value = 1728.1034
$ws.Range("I126").Value = 1757.5555
$ws.Range("J126").Value = 1679.909
$ws.Range("K126").Value = 5272.666499999999
$ws.Range("L126").Value = 5039.727000000001
$ws.Range("M126").Value = -2802.666499999999
$ws.Range("N126").Value = -9979.727000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 776.4
$ws.Range("I5").Value = 609.1429000000001
$ws.Range("J5").Value = 1166.6666
$ws.Range("K5").Value = 1827.4287
$ws.Range("L5").Value = 3499.9998
$ws.Range("M5").Value = -1715.4287
$ws.Range("N5").Value = -3723.9998

$ws.Range("H117").Value = 3577.4092
$ws.Range("J117").Value = 4486.1177
$ws.Range("L117").Value = 13458.3531
$ws.Range("N117").Value = -20342.3531

$ws.Range("H129").Value = 53928.95
$ws.Range("I129").Value = 72133.57000000001
$ws.Range("K129").Value = 216400.71
$ws.Range("M129").Value = -211400.71

$ws.Range("H135").Value = 776.4
$ws.Range("I135").Value = 609.1429000000001
$ws.Range("J135").Value = 1166.6666
$ws.Range("K135").Value = 5482.2861
$ws.Range("L135").Value = 10499.9994
$ws.Range("M135").Value = -2947.2861
$ws.Range("N135").Value = -15569.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2601.6
$ws.Range("I97").Value = 2752.25
$ws.Range("J97").Value = 1999
$ws.Range("K97").Value = 2752.25
$ws.Range("L97").Value = 1999
$ws.Range("M97").Value = -2256.25
$ws.Range("N97").Value = -2991

$ws.Range("H102").Value = 2193.818
$ws.Range("I102").Value = 1849
$ws.Range("J102").Value = 3113.3333
$ws.Range("K102").Value = 1849
$ws.Range("L102").Value = 3113.3333
$ws.Range("M102").Value = -227
$ws.Range("N102").Value = -6357.3333

$ws.Range("H122").Value = 2300
$ws.Range("I122").Value = 2120
$ws.Range("J122").Value = 2750
$ws.Range("K122").Value = 6360
$ws.Range("L122").Value = 8250
$ws.Range("M122").Value = -3910
$ws.Range("N122").Value = -13150

$ws.Range("H123").Value = 38450
$ws.Range("J123").Value = 38450
$ws.Range("L123").Value = 38450
$ws.Range("N123").Value = -43350

$ws.Range("H132").Value = 3413.5
$ws.Range("I132").Value = 3300
$ws.Range("J132").Value = 3527
$ws.Range("K132").Value = 9900
$ws.Range("L132").Value = 10581
$ws.Range("M132").Value = -7370
$ws.Range("N132").Value = -15641

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3404.7222
$ws.Range("I122").Value = 2966.6667
$ws.Range("J122").Value = 3623.75
$ws.Range("K122").Value = 8900.000100000001
$ws.Range("L122").Value = 10871.25
$ws.Range("M122").Value = -6450.000100000001
$ws.Range("N122").Value = -15771.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3900
$ws.Range("I62").Value = 4250
$ws.Range("J62").Value = 3783.3333
$ws.Range("K62").Value = 4250
$ws.Range("L62").Value = 3783.3333
$ws.Range("M62").Value = -3626
$ws.Range("N62").Value = -5031.3333

$ws.Range("H65").Value = 3900
$ws.Range("I65").Value = 4250
$ws.Range("J65").Value = 3783.3333
$ws.Range("K65").Value = 21250
$ws.Range("L65").Value = 18916.6665
$ws.Range("M65").Value = -18130
$ws.Range("N65").Value = -25156.6665

$ws.Range("H109").Value = 15662.333
$ws.Range("J109").Value = 15662.333
$ws.Range("L109").Value = 15662.333
$ws.Range("N109").Value = -18436.333

$ws.Range("H122").Value = 3183.0527
$ws.Range("I122").Value = 2607.8
$ws.Range("J122").Value = 3822.2222
$ws.Range("K122").Value = 7823.400000000001
$ws.Range("L122").Value = 11466.6666
$ws.Range("M122").Value = -5373.400000000001
$ws.Range("N122").Value = -16366.6666

$ws.Range("H126").Value = 2115.875
$ws.Range("I126").Value = 1464.4286
$ws.Range("J126").Value = 3027.9
$ws.Range("K126").Value = 4393.2858
$ws.Range("L126").Value = 9083.700000000001
$ws.Range("M126").Value = -1923.2858
$ws.Range("N126").Value = -14023.7
